$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell content updates ---
# Row 2: switch the extension-script URL to the raw.githubusercontent.com host
# and bump the VM size from Basic_A0 to Standard_DS1_v2
$ws.Range("C2").Value = "https://raw.githubusercontent.com/ansible/ansible/devel/examples/scripts/ConfigureRemotingForAnsible.ps1"
$ws.Range("D2").Value = "Standard_DS1_v2"

# Row 3: use "Windows Server 2016 Datacenter" (since Core isn't available), same URL + VM size bump
$ws.Range("B3").Value = "Windows Server 2016 Datacenter"
$ws.Range("C3").Value = "https://raw.githubusercontent.com/ansible/ansible/devel/examples/scripts/ConfigureRemotingForAnsible.ps1"
$ws.Range("D3").Value = "Standard_DS1_v2"

# Row 5/6: tidy up the OS labels
$ws.Range("B5").Value = "RHEL 7"
$ws.Range("B6").Value = "Ubuntu 14"

# --- Formatting: make B5/B6 bold like the other App Name header cells (A2/B2/B3) ---
$ws.Range("A2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column width tweaks ---
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(7).EntireColumn.AutoFit()

# --- Selection moves to D3 ---
$ws.Range("D3").Select()
